$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 6
$ws.Range("AI2").Value = 13
$ws.Range("AW2").Value = 4.75
$ws.Range("G2").Value = 2.6
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3.4
$ws.Range("L2").Value = 3.75
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 26
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.36
$ws.Range("M4").Value = 1.04
$ws.Range("O4").Value = 1.22
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 451
$ws.Range("AH5").Value = 17
$ws.Range("AL5").Value = 51
$ws.Range("AN5").Value = 3.25
$ws.Range("AU5").Value = 9.5
$ws.Range("AW5").Value = 8.5
$ws.Range("AZ5").Value = 151
$ws.Range("G5").Value = 1.4
$ws.Range("H5").Value = 4.75
$ws.Range("I5").Value = 8
$ws.Range("J5").Value = 1.91
$ws.Range("K5").Value = 2.38
$ws.Range("L5").Value = 7.5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("Q5").Value = 1.87
$ws.Range("R5").Value = 2.03
$ws.Range("U5").Value = 2.05
$ws.Range("V5").Value = 1.7
$ws.Range("X5").Value = 6.5
$ws.Range("Y5").Value = 8.5
$ws.Range("Z5").Value = 9
$ws.Range("AH6").Value = 13
$ws.Range("AT6").Value = 2.63
$ws.Range("AV6").Value = 67
$ws.Range("AW6").Value = 7.5
$ws.Range("G6").Value = 1.6
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 6.25
$ws.Range("J6").Value = 2.25
$ws.Range("S6").Value = 1.44
$ws.Range("T6").Value = 2.63
$ws.Range("O7").Value = 1.29
$ws.Range("P7").Value = 3.5
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.9
